$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as plain text, preserving the original "General"
# number format/style (avoids Excel auto-converting numeric-looking
# strings such as "597.18" or "65.110.26" into floating point numbers).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Row 29/30 swap: PEPE moves to row 29, Binance-PegBSC-USD moves to row 30 ---
$sub3 = [char]0x2083
$pepePrice = "{0}{1}{2}" -f "0.0", $sub3, "0952"

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D29") $pepePrice
$ws.Range("E29").Value = "  +10.73%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  +0.03%  "

# --- Price (D) and Volume(1h) (E) updates for remaining rows ---
Set-TextValue $ws.Range("D2") "65.110.26"
$ws.Range("E2").Value = "  +2.74%  "
Set-TextValue $ws.Range("D3") "2.637.23"
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue $ws.Range("D5") "597.18"
$ws.Range("E5").Value = "  +1.71%  "
Set-TextValue $ws.Range("D6") "155.96"
$ws.Range("E6").Value = "  +4.58%  "
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.09%  "
Set-TextValue $ws.Range("D8") "0.592"
$ws.Range("E8").Value = "  +1.11%  "
Set-TextValue $ws.Range("D9") "0.118"
$ws.Range("E9").Value = "  +7.72%  "
$ws.Range("E10").Value = "  +4.94%  "
$ws.Range("E11").Value = "  +0.55%  "
Set-TextValue $ws.Range("D12") "0.153"
$ws.Range("E12").Value = "  +2.31%  "
Set-TextValue $ws.Range("D13") "29.29"
$ws.Range("E13").Value = "  +6.63%  "
Set-TextValue $ws.Range("D14") "0.0000188"
$ws.Range("E14").Value = "  +21.83%  "
Set-TextValue $ws.Range("D15") "3.110.09"
$ws.Range("E15").Value = "  +2.31%  "
Set-TextValue $ws.Range("D16") "64.998.89"
$ws.Range("E16").Value = "  +2.79%  "
Set-TextValue $ws.Range("D17") "2.644.65"
$ws.Range("E17").Value = "  +2.39%  "
Set-TextValue $ws.Range("D18") "12.57"
$ws.Range("E18").Value = "  +3.38%  "
Set-TextValue $ws.Range("D19") "4.81"
$ws.Range("E19").Value = "  +2.96%  "
Set-TextValue $ws.Range("D20") "353.03"
$ws.Range("E20").Value = "  +2.30%  "
Set-TextValue $ws.Range("D21") "7.35"
$ws.Range("E21").Value = "  +7.73%  "
$ws.Range("E22").Value = "  +0.16%  "
Set-TextValue $ws.Range("D23") "68.34"
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("E24").Value = "  +0.42%  "
Set-TextValue $ws.Range("D25") "9.55"
$ws.Range("E25").Value = "  +4.95%  "
$ws.Range("E26").Value = "  -1.02%  "
$ws.Range("E27").Value = "  +1.11%  "
Set-TextValue $ws.Range("D28") "8.09"
$ws.Range("E28").Value = "  +0.97%  "
Set-TextValue $ws.Range("D31") "2.11"
$ws.Range("E31").Value = "  +4.36%  "
Set-TextValue $ws.Range("D32") "510.18"
$ws.Range("E32").Value = "  -7.78%  "
Set-TextValue $ws.Range("D33") "1.77"
$ws.Range("E33").Value = "  +1.40%  "
Set-TextValue $ws.Range("D34") "5.61"
$ws.Range("E34").Value = "  +8.09%  "
Set-TextValue $ws.Range("D35") "6.35"
$ws.Range("E35").Value = "  +6.07%  "
$ws.Range("E36").Value = "  +2.80%  "
Set-TextValue $ws.Range("D37") "20.29"
$ws.Range("E37").Value = "  +4.60%  "
Set-TextValue $ws.Range("D38") "163.51"
$ws.Range("E38").Value = "  -1.09%  "
Set-TextValue $ws.Range("D39") "2.01"
$ws.Range("E39").Value = "  +6.13%  "
Set-TextValue $ws.Range("D40") "0.998"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("E41").Value = "  -0.01%  "
Set-TextValue $ws.Range("D42") "42.29"
$ws.Range("E42").Value = "  +6.47%  "
Set-TextValue $ws.Range("D43") "165.60"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  +2.92%  "
$ws.Range("E45").Value = "  +4.30%  "
Set-TextValue $ws.Range("D46") "23.16"
$ws.Range("E46").Value = "  +2.03%  "
Set-TextValue $ws.Range("D47") "2.21"
$ws.Range("E47").Value = "  +8.47%  "
Set-TextValue $ws.Range("D48") "0.647"
$ws.Range("E48").Value = "  +3.37%  "
$ws.Range("E49").Value = "  +2.73%  "
Set-TextValue $ws.Range("D50") "0.0984"
$ws.Range("E50").Value = "  +2.53%  "
Set-TextValue $ws.Range("D51") "19.44"
$ws.Range("E51").Value = "  +2.52%  "
